$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 18 and 19 (player "Ja Morant" / "De'Aaron Fox") need their
# Player Name (A) and Team (C) values swapped. Position (B) is "PG"
# for both rows already, so it is left untouched.

$a18 = $ws.Range("A18").Value2
$c18 = $ws.Range("C18").Value2
$a19 = $ws.Range("A19").Value2
$c19 = $ws.Range("C19").Value2

$ws.Range("A18").Value2 = $a19
$ws.Range("C18").Value2 = $c19
$ws.Range("A19").Value2 = $a18
$ws.Range("C19").Value2 = $c18
